$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook (after the current last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "checkArtifactsManager"

# Populate the two cells (these become new shared-string entries).
$newSheet.Range("A1").Value = "artifactsManagerTitle"
$newSheet.Range("A2").Value = "Artifacts Manager"

# Match the column width used elsewhere (raw OOXML width 27.5).
$newSheet.Columns.Item(1).ColumnWidth = 26.6666666666667

# Make the new sheet the active / selected tab, with A2 selected.
$newSheet.Range("A2").Select()
$newSheet.Activate()
